# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# All Coin/Link/Price/Volume(1h) cells in this sheet are stored as literal TEXT
# (OOXML inlineStr), not numbers -- Price strings use "." as both a thousands
# separator and a decimal point (e.g. "69.044.46"), which would be corrupted if
# treated as a number. Plain `Range.Value = "123.45"` lets Excel auto-detect a
# real number for any unambiguous decimal, which can (a) drop meaningful trailing
# zeros ("38.00" -> 38) and (b) round-trip through a binary double and pick up
# float noise ("601.91" -> 601.90999999999997) even though Value2/Text still
# *display* the clean string. Set-TextValue forces text storage for every Price
# write via a throwaway "@" (text) number format, then clears the format again so
# the cell keeps its original (unstyled) look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 'D2' '69.019.65'
$ws.Range('E2').Value = '  +0.58%  '

Set-TextValue 'D3' '3.739.38'
$ws.Range('E3').Value = '  +1.14%  '

Set-TextValue 'D5' '601.91'
$ws.Range('E5').Value = '  +0.46%  '

Set-TextValue 'D6' '167.64'
$ws.Range('E6').Value = '  +0.44%  '

Set-TextValue 'D7' '3.736.69'
$ws.Range('E7').Value = '  +1.00%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('E9').Value = '  +0.87%  '

$ws.Range('E10').Value = '  +2.78%  '

$ws.Range('E11').Value = '  +2.91%  '

$ws.Range('E12').Value = '  +0.48%  '

Set-TextValue 'D13' '38.00'
$ws.Range('E13').Value = '  +0.03%  '

$ws.Range('E14').Value = '  +2.22%  '

Set-TextValue 'D15' '4.368.15'
$ws.Range('E15').Value = '  +1.26%  '

Set-TextValue 'D16' '3.740.31'
$ws.Range('E16').Value = '  +1.20%  '

Set-TextValue 'D17' '69.027.95'
$ws.Range('E17').Value = '  +0.58%  '

$ws.Range('E18').Value = '  +1.43%  '

Set-TextValue 'D19' '17.29'
$ws.Range('E19').Value = '  +1.05%  '

$ws.Range('E20').Value = '  -1.30%  '

Set-TextValue 'D21' '10.94'
$ws.Range('E21').Value = '  +19.96%  '

Set-TextValue 'D22' '492.67'
$ws.Range('E22').Value = '  +0.07%  '

Set-TextValue 'D23' '0.726'
$ws.Range('E23').Value = '  +0.80%  '

$ws.Range('E24').Value = '  +8.69%  '

Set-TextValue 'D25' '84.80'
$ws.Range('E25').Value = '  +0.53%  '

$ws.Range('E26').Value = '  +0.73%  '

Set-TextValue 'D27' '12.34'
$ws.Range('E27').Value = '  +1.46%  '

Set-TextValue 'D28' '10.10'
$ws.Range('E28').Value = '  +0.72%  '

$ws.Range('E29').Value = '  -0.10%  '

$ws.Range('E30').Value = '  +2.49%  '

Set-TextValue 'D31' '2.47'
$ws.Range('E31').Value = '  +4.57%  '

$ws.Range('E32').Value = '  +3.27%  '

Set-TextValue 'D33' '31.49'
$ws.Range('E33').Value = '  +0.28%  '

Set-TextValue 'D34' '3.885.97'
$ws.Range('E34').Value = '  +1.42%  '

$ws.Range('E35').Value = '  +0.43%  '

Set-TextValue 'D36' '3.677.58'
$ws.Range('E36').Value = '  +1.15%  '

Set-TextValue 'D37' '1.00'
$ws.Range('E37').Value = '  -0.02%  '

$ws.Range('E38').Value = '  +1.86%  '

$ws.Range('E39').Value = '  +3.64%  '

$ws.Range('E40').Value = '  +1.81%  '

Set-TextValue 'D41' '0.324'
$ws.Range('E41').Value = '  +0.80%  '

$ws.Range('E42').Value = '  +6.25%  '

Set-TextValue 'D43' '432.13'
$ws.Range('E43').Value = '  +0.30%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D44' '48.65'
$ws.Range('E44').Value = '  -0.87%  '

$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D45' '1.99'
$ws.Range('E45').Value = '  +1.46%  '

$ws.Range('E46').Value = '  +1.39%  '

$ws.Range('E47').Value = '  +0.00%  '

Set-TextValue 'D48' '40.34'
$ws.Range('E48').Value = '  +0.50%  '

Set-TextValue 'D49' '141.15'
$ws.Range('E49').Value = '  -0.35%  '

Set-TextValue 'D50' '2.782.83'
$ws.Range('E50').Value = '  +1.90%  '

$ws.Range('E51').Value = '  +1.01%  '
